$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: tortellini -> ragu
$ws.Range("A2").Value = "ragu"
$ws.Range("B2").Value = 0.6830000000000001
$ws.Range("C2").Value = 6.83
$ws.Range("D2").Value = 0.6830000000000001
$ws.Range("E2").Value = '{"name": "food (food)", "total": {"name": [], "protein": 0.515, "calories": 5.15, "cost": 0.515, "amount": 0.0}}'

# Row 3: tortellini -> ragu
$ws.Range("A3").Value = "ragu"
$ws.Range("B3").Value = 0.6830000000000001
$ws.Range("C3").Value = 6.83
$ws.Range("D3").Value = 0.6830000000000001
$ws.Range("E3").Value = '{"name": "food name (food vendor)", "total": {"name": [], "protein": 0.168, "calories": 1.68, "cost": 0.168, "amount": 0.0}}'
